$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

function Set-TextValue($cell, $value) {
    # Force a numeric-looking string to be stored as text (shared string)
    # instead of being auto-coerced to a number, while leaving the cell's
    # style back at the default (no leftover explicit number format).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-HyperlinkEmail($cell, $address) {
    $cell.Value = $address
    if ($cell.Hyperlinks.Count -gt 0) { $cell.Hyperlinks.Delete() }
    $ws.Hyperlinks.Add($cell, "mailto:" + $address) | Out-Null
    $cell.Style = "Hyperlink"
}

# Duplicate row 2 into a brand new row 3, preserving cell styles
# (hyperlink font, numeric formatting, etc.) via a row copy + insert.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(2, 8).Value = 1234567890
$ws.Cells.Item(3, 8).Value = 1234567890
$ws.Cells.Item(3, 2).Value = "f"

Set-HyperlinkEmail $ws.Cells.Item(2, 9) "testacc3185@gmail.com"
Set-HyperlinkEmail $ws.Cells.Item(3, 9) "testacc3186@gmail.com"

Set-TextValue $ws.Cells.Item(2, 11) "57870"
Set-TextValue $ws.Cells.Item(3, 11) "84327"

Set-TextValue $ws.Cells.Item(2, 14) "58258"
Set-TextValue $ws.Cells.Item(3, 14) "58259"

# Move the active selection to I4, matching the saved view state.
$ws.Cells.Item(4, 9).Select()
